$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the order of items in row 2 (A2): "2-queque,1-torta," -> "1-torta,2-queque,"
$ws.Range("A2").Value = "1-torta,2-queque,"

# Add new order row 9 - force text format so values like dates and long
# digit strings are stored as plain text (matching columns A-G elsewhere).
$ws.Range("A9:G9").NumberFormat = "@"
$ws.Range("A9").Value = "1-torta,"
$ws.Range("B9").Value = "10-08-2020"
$ws.Range("C9").Value = "12-10-2019"
$ws.Range("D9").Value = "lala alal "
$ws.Range("E9").Value = "ii@ii.com"
$ws.Range("F9").Value = "999999999"
$ws.Range("G9").Value = "En Proceso"
$ws.Range("H9").Value = 30000
$ws.Range("I9").Value = 16000
$ws.Range("J9").Value = 0
